$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B / C (coin name / link) updates ---
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

# --- Column E (Volume/1h %) updates ---
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('E8').Value = '  -1.27%  '
$ws.Range('E9').Value = '  +8.61%  '
$ws.Range('E10').Value = '  +1.43%  '
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('E15').Value = '  -0.93%  '
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('E18').Value = '  +3.81%  '
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('E20').Value = '  +0.71%  '
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('E22').Value = '  -2.20%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('E24').Value = '  +3.35%  '
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('E26').Value = '  +2.25%  '
$ws.Range('E27').Value = '  +1.84%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  +0.29%  '
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('E31').Value = '  +2.09%  '
$ws.Range('E32').Value = '  -2.55%  '
$ws.Range('E33').Value = '  +1.06%  '
$ws.Range('E34').Value = '  +7.35%  '
$ws.Range('E35').Value = '  +0.80%  '
$ws.Range('E36').Value = '  +9.20%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('E38').Value = '  +1.74%  '
$ws.Range('E39').Value = '  +2.25%  '
$ws.Range('E40').Value = '  -2.27%  '
$ws.Range('E41').Value = '  +1.79%  '
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('E44').Value = '  +5.36%  '
$ws.Range('E45').Value = '  +0.87%  '
$ws.Range('E46').Value = '  +2.48%  '
$ws.Range('E47').Value = '  +2.76%  '
$ws.Range('E48').Value = '  +1.66%  '
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('E50').Value = '  +0.45%  '

# --- Column D (Price) updates: force text to avoid numeric coercion ---
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.484.74'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.825.94'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5172'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3859'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08297'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.124'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.94'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.376'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.474'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.823.40'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001123'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06637'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.050'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.511.37'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.51'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.249'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.10'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.61'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.034.45'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.01'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1109'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.092'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.735'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07536'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.688'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.26'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2225'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02364'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.251'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.758'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6401'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.185'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.397'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6218'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.804'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '127.69'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.011'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.205'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06960'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.080'
$ws.Range('D51').Style = 'Normal'
